$d = $word.ActiveDocument

# Update the date/title line (unique text, safe to Find & Replace across the document)
$d.Content.Find.Execute("2026-01-01 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-02 Friday", 2) | Out-Null

# Update each multiplication-fact cell directly by table coordinates.
# (Several cells share old/new text values across rows, e.g. "432×6=2592" is both a
#  target of one cell and the prior value of another, so a document-wide Find/Replace
#  could cross-contaminate; addressing cells directly avoids that.)
$tbl = $d.Tables.Item(1)

$cell = $tbl.Cell(1, 1)
if ($cell.Range.Text -notlike "274×7=1918*") { throw "Unexpected text in cell(1,1): $($cell.Range.Text)" }
$cell.Range.Text = "867×8=6936"

$cell = $tbl.Cell(1, 2)
if ($cell.Range.Text -notlike "145×2=290*") { throw "Unexpected text in cell(1,2): $($cell.Range.Text)" }
$cell.Range.Text = "703×7=4921"

$cell = $tbl.Cell(1, 3)
if ($cell.Range.Text -notlike "654×5=3270*") { throw "Unexpected text in cell(1,3): $($cell.Range.Text)" }
$cell.Range.Text = "749×6=4494"

$cell = $tbl.Cell(1, 4)
if ($cell.Range.Text -notlike "555×4=2220*") { throw "Unexpected text in cell(1,4): $($cell.Range.Text)" }
$cell.Range.Text = "432×6=2592"

$cell = $tbl.Cell(1, 5)
if ($cell.Range.Text -notlike "531×7=3717*") { throw "Unexpected text in cell(1,5): $($cell.Range.Text)" }
$cell.Range.Text = "111×7=777"

$cell = $tbl.Cell(5, 1)
if ($cell.Range.Text -notlike "621×7=4347*") { throw "Unexpected text in cell(5,1): $($cell.Range.Text)" }
$cell.Range.Text = "401×4=1604"

$cell = $tbl.Cell(5, 2)
if ($cell.Range.Text -notlike "791×4=3164*") { throw "Unexpected text in cell(5,2): $($cell.Range.Text)" }
$cell.Range.Text = "209×6=1254"

$cell = $tbl.Cell(5, 3)
if ($cell.Range.Text -notlike "101×2=202*") { throw "Unexpected text in cell(5,3): $($cell.Range.Text)" }
$cell.Range.Text = "643×5=3215"

$cell = $tbl.Cell(5, 4)
if ($cell.Range.Text -notlike "415×5=2075*") { throw "Unexpected text in cell(5,4): $($cell.Range.Text)" }
$cell.Range.Text = "661×7=4627"

$cell = $tbl.Cell(5, 5)
if ($cell.Range.Text -notlike "640×6=3840*") { throw "Unexpected text in cell(5,5): $($cell.Range.Text)" }
$cell.Range.Text = "150×2=300"

$cell = $tbl.Cell(10, 1)
if ($cell.Range.Text -notlike "302×2=604*") { throw "Unexpected text in cell(10,1): $($cell.Range.Text)" }
$cell.Range.Text = "439×8=3512"

$cell = $tbl.Cell(10, 2)
if ($cell.Range.Text -notlike "493×4=1972*") { throw "Unexpected text in cell(10,2): $($cell.Range.Text)" }
$cell.Range.Text = "184×2=368"

$cell = $tbl.Cell(10, 3)
if ($cell.Range.Text -notlike "447×9=4023*") { throw "Unexpected text in cell(10,3): $($cell.Range.Text)" }
$cell.Range.Text = "443×2=886"

$cell = $tbl.Cell(10, 4)
if ($cell.Range.Text -notlike "112×7=784*") { throw "Unexpected text in cell(10,4): $($cell.Range.Text)" }
$cell.Range.Text = "309×5=1545"

$cell = $tbl.Cell(10, 5)
if ($cell.Range.Text -notlike "627×8=5016*") { throw "Unexpected text in cell(10,5): $($cell.Range.Text)" }
$cell.Range.Text = "805×4=3220"

$cell = $tbl.Cell(15, 1)
if ($cell.Range.Text -notlike "902×6=5412*") { throw "Unexpected text in cell(15,1): $($cell.Range.Text)" }
$cell.Range.Text = "905×6=5430"

$cell = $tbl.Cell(15, 2)
if ($cell.Range.Text -notlike "367×8=2936*") { throw "Unexpected text in cell(15,2): $($cell.Range.Text)" }
$cell.Range.Text = "549×7=3843"

$cell = $tbl.Cell(15, 3)
if ($cell.Range.Text -notlike "190×7=1330*") { throw "Unexpected text in cell(15,3): $($cell.Range.Text)" }
$cell.Range.Text = "476×3=1428"

$cell = $tbl.Cell(15, 4)
if ($cell.Range.Text -notlike "432×6=2592*") { throw "Unexpected text in cell(15,4): $($cell.Range.Text)" }
$cell.Range.Text = "435×4=1740"

$cell = $tbl.Cell(15, 5)
if ($cell.Range.Text -notlike "143×8=1144*") { throw "Unexpected text in cell(15,5): $($cell.Range.Text)" }
$cell.Range.Text = "255×3=765"

$cell = $tbl.Cell(20, 1)
if ($cell.Range.Text -notlike "565×4=2260*") { throw "Unexpected text in cell(20,1): $($cell.Range.Text)" }
$cell.Range.Text = "890×4=3560"

$cell = $tbl.Cell(20, 2)
if ($cell.Range.Text -notlike "476×4=1904*") { throw "Unexpected text in cell(20,2): $($cell.Range.Text)" }
$cell.Range.Text = "872×7=6104"

$cell = $tbl.Cell(20, 3)
if ($cell.Range.Text -notlike "556×4=2224*") { throw "Unexpected text in cell(20,3): $($cell.Range.Text)" }
$cell.Range.Text = "804×9=7236"

$cell = $tbl.Cell(20, 4)
if ($cell.Range.Text -notlike "231×9=2079*") { throw "Unexpected text in cell(20,4): $($cell.Range.Text)" }
$cell.Range.Text = "574×9=5166"

$cell = $tbl.Cell(20, 5)
if ($cell.Range.Text -notlike "137×3=411*") { throw "Unexpected text in cell(20,5): $($cell.Range.Text)" }
$cell.Range.Text = "322×4=1288"
